$d = $word.ActiveDocument

$d.Content.Find.Execute("1 Tim. 5:17-25", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1 Tim. 5:17-6:2", 2) | Out-Null

$d.Content.Find.Execute("1 Tim. 6:1-21", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1 Tim. 6:3-21", 2) | Out-Null
